$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The account-statement database was refreshed: the two oldest workers on the
# employee table (EVA MARGARITA ANGULO MERCADO / 22789921 and ERIKA MILENA
# RAMOS BARRERA / 45530409, rows 16-17, period 2507) are removed. Deleting
# those two rows shifts the remaining two workers (STELLA ANGULO CARMONA and
# TIVEL ESTEPHANI BATISTA MARTINEZ) up into rows 16-17, and moves the
# signature block up from rows 24-25 to rows 22-23.
$ws.Rows("16:17").Delete()

# The remaining two workers now belong to the new "Periodo Mora" 2508
# (previously 2507).
$ws.Range("E16").Value = "2508"
$ws.Range("E17").Value = "2508"

# Update the summary totals to reflect the new (smaller) worker list.
$ws.Range("E11").Value = 260000
$ws.Range("C13").Value = 2
